$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G7").Value = "make the servo pull string faster"
$ws.Range("G8").Value = "set up puppet with strings and servos"
$ws.Range("G9").Value = "test with software"
$ws.Range("G10").Value = "build frame for puppet display"
$ws.Range("G11").Value = "experiment with motor to move entire puppet left and right"
$ws.Range("G12").Value = "build frame rail for puppet movement"
$ws.Range("G13").Value = "create movement segments"
$ws.Range("G14").Value = "prepare and test plan for disassembly and reassembly of hardware"
$ws.Range("G15").Value = ""
$ws.Range("G16").Value = ""
$ws.Range("G17").Value = ""

$null = $ws.Range("G7").Select()
